{"js": "const pairs = [\n  [\"2024-05-18 Saturday\", \"2024-05-19 Sunday\"],\n  [\"826\u00d74=\", \"956\u00d75=\"],\n  [\"470\u00d77=\", \"216\u00d73=\"],\n  [\"292\u00d75=\", \"615\u00d78=\"],\n  [\"657\u00d78=\", \"356\u00d74=\"],\n  [\"777\u00d77=\", \"258\u00d75=\"],\n  [\"253\u00d79=\", \"171\u00d73=\"],\n  [\"731\u00d78=\", \"340\u00d75=\"],\n  [\"825\u00d74=\", \"326\u00d74=\"],\n  [\"186\u00d73=\", \"129\u00d73=\"],\n  [\"631\u00d77=\", \"658\u00d75=\"],\n  [\"341\u00d78=\", \"593\u00d73=\"],\n  [\"586\u00d75=\", \"244\u00d79=\"],\n  [\"181\u00d73=\", \"520\u00d73=\"],\n  [\"209\u00d76=\", \"952\u00d73=\"],\n  [\"885\u00d74=\", \"549\u00d79=\"],\n  [\"115\u00d72=\", \"880\u00d79=\"],\n  [\"555\u00d75=\", \"168\u00d78=\"],\n  [\"824\u00d74=\", \"155\u00d73=\"],\n  [\"403\u00d76=\", \"797\u00d75=\"],\n  [\"651\u00d77=\", \"791\u00d79=\"],\n  [\"276\u00d78=\", \"202\u00d78=\"],\n  [\"772\u00d77=\", \"464\u00d74=\"],\n  [\"664\u00d79=\", \"255\u00d75=\"],\n  [\"323\u00d75=\", \"919\u00d78=\"],\n  [\"929\u00d78=\", \"343\u00d74=\"],\n];\n\nfor (const [before, after] of pairs) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + before);\n  }\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-05-18 Saturday\", \"2024-05-19 Sunday\"),\n    @(\"826\u00d74=\", \"956\u00d75=\"),\n    @(\"470\u00d77=\", \"216\u00d73=\"),\n    @(\"292\u00d75=\", \"615\u00d78=\"),\n    @(\"657\u00d78=\", \"356\u00d74=\"),\n    @(\"777\u00d77=\", \"258\u00d75=\"),\n    @(\"253\u00d79=\", \"171\u00d73=\"),\n    @(\"731\u00d78=\", \"340\u00d75=\"),\n    @(\"825\u00d74=\", \"326\u00d74=\"),\n    @(\"186\u00d73=\", \"129\u00d73=\"),\n    @(\"631\u00d77=\", \"658\u00d75=\"),\n    @(\"341\u00d78=\", \"593\u00d73=\"),\n    @(\"586\u00d75=\", \"244\u00d79=\"),\n    @(\"181\u00d73=\", \"520\u00d73=\"),\n    @(\"209\u00d76=\", \"952\u00d73=\"),\n    @(\"885\u00d74=\", \"549\u00d79=\"),\n    @(\"115\u00d72=\", \"880\u00d79=\"),\n    @(\"555\u00d75=\", \"168\u00d78=\"),\n    @(\"824\u00d74=\", \"155\u00d73=\"),\n    @(\"403\u00d76=\", \"797\u00d75=\"),\n    @(\"651\u00d77=\", \"791\u00d79=\"),\n    @(\"276\u00d78=\", \"202\u00d78=\"),\n    @(\"772\u00d77=\", \"464\u00d74=\"),\n    @(\"664\u00d79=\", \"255\u00d75=\"),\n    @(\"323\u00d75=\", \"919\u00d78=\"),\n    @(\"929\u00d78=\", \"343\u00d74=\"),\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($pair[0], $false, $true, $false, $false, $false, $true, 0, $false, $pair[1], 2) | Out-Null\n}\n"}
